# Update the contacts table: replace the LastName and email columns
# with the new values (firstName/Company stay the same), and move the
# active selection to E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LastName column (B)
$ws.Range("B2").Value = "Dia"
$ws.Range("B3").Value = "Dio"
$ws.Range("B4").Value = "Dpo"

# email column (D)
$ws.Range("D2").Value = "joee@yopmail.com"
$ws.Range("D3").Value = "adile@yopmail.com"
$ws.Range("D4").Value = "nabil@yopmail.com"

# Move the selection to E13, matching the saved view state.
$ws.Range("E13").Select()
